$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fully clear former data region (A2:D44) and extend through new rows (A2:D47)
$ws.Range("A2:D47").Clear()

$ws.Range("B2").Value = "TV"
$ws.Range("C2").Value = "Spend"
$ws.Range("D2").Value = 147

$ws.Range("A3").Value = 46184
$ws.Range("A3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B3").Value = "TV"
$ws.Range("C3").Value = "Spend"
$ws.Range("D3").Value = 83

$ws.Range("A4").Value = 46226
$ws.Range("A4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B4").Value = "Radio"
$ws.Range("C4").Value = "Spend"
$ws.Range("D4").Value = 68

$ws.Range("A5").Value = 46219
$ws.Range("A5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B5").Value = "Radio"
$ws.Range("C5").Value = "Spend"

$ws.Range("A6").Value = 46205
$ws.Range("A6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D6").Value = 95

$ws.Range("A7").Value = 46219
$ws.Range("A7").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C7").Value = "GRPs"
$ws.Range("D7").Value = 10

$ws.Range("A8").Value = 46198
$ws.Range("A8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C8").Value = "Spend"
$ws.Range("D8").Value = 166

$ws.Range("B9").Value = "Radio"
$ws.Range("C9").Value = "Spend"
$ws.Range("D9").Value = 177

$ws.Range("A10").Value = 46191
$ws.Range("A10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B10").Value = "TV"
$ws.Range("C10").Value = "Spend"

$ws.Range("A11").Value = 46226
$ws.Range("A11").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B11").Value = "TV"
$ws.Range("C11").Value = "Spend"

$ws.Range("B12").Value = "TV"
$ws.Range("D12").Value = 1

$ws.Range("A13").Value = 46219
$ws.Range("A13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B13").Value = "TV"
$ws.Range("C13").Value = "GRPs"
$ws.Range("D13").Value = 6

$ws.Range("A14").Value = 46198
$ws.Range("A14").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B14").Value = "TV"
$ws.Range("D14").Value = 9

$ws.Range("A15").Value = 46247
$ws.Range("A15").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B15").Value = "TV"
$ws.Range("D15").Value = 151

$ws.Range("A16").Value = 46254
$ws.Range("A16").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B16").Value = "TV"
$ws.Range("C16").Value = "Spend"
$ws.Range("D16").Value = 139

$ws.Range("A17").Value = 46219
$ws.Range("A17").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B17").Value = "TV"
$ws.Range("C17").Value = "Spend"
$ws.Range("D17").Value = 174

$ws.Range("B18").Value = "TV"
$ws.Range("C18").Value = "GRPs"
$ws.Range("D18").Value = 6

$ws.Range("A19").Value = 46240
$ws.Range("A19").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B19").Value = "Radio"
$ws.Range("C19").Value = "GRPs"
$ws.Range("D19").Value = 8

$ws.Range("A20").Value = 46198
$ws.Range("A20").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B20").Value = "Radio"
$ws.Range("C20").Value = "GRPs"

$ws.Range("B21").Value = "TV"
$ws.Range("D21").Value = 10

$ws.Range("A22").Value = 46198
$ws.Range("A22").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B22").Value = "Radio"
$ws.Range("C22").Value = "GRPs"
$ws.Range("D22").Value = 10

$ws.Range("A23").Value = 46233
$ws.Range("A23").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B23").Value = "Radio"
$ws.Range("C23").Value = "Spend"
$ws.Range("D23").Value = 135

$ws.Range("A24").Value = 46254
$ws.Range("A24").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B24").Value = "TV"
$ws.Range("C24").Value = "Spend"
$ws.Range("D24").Value = 139

$ws.Range("A25").Value = 46233
$ws.Range("A25").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C25").Value = "Spend"

$ws.Range("A26").Value = 46198
$ws.Range("A26").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B26").Value = "TV"
$ws.Range("C26").Value = "GRPs"
$ws.Range("D26").Value = 9

$ws.Range("A27").Value = 46219
$ws.Range("A27").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B27").Value = "Radio"
$ws.Range("D27").Value = 121

$ws.Range("A28").Value = 46191
$ws.Range("A28").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B28").Value = "TV"
$ws.Range("C28").Value = "GRPs"
$ws.Range("D28").Value = 6

$ws.Range("A29").Value = 46247
$ws.Range("A29").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B29").Value = "Radio"
$ws.Range("C29").Value = "GRPs"

$ws.Range("A30").Value = 46233
$ws.Range("A30").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B30").Value = "Radio"
$ws.Range("C30").Value = "GRPs"

$ws.Range("A31").Value = 46240
$ws.Range("A31").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B31").Value = "Radio"
$ws.Range("C31").Value = "GRPs"

$ws.Range("A32").Value = 46184
$ws.Range("A32").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C32").Value = "GRPs"
$ws.Range("D32").Value = 4

$ws.Range("A33").Value = 46184
$ws.Range("A33").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B33").Value = "TV"
$ws.Range("C33").Value = "Spend"
$ws.Range("D33").Value = 83

$ws.Range("A34").Value = 46212
$ws.Range("A34").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B34").Value = "TV"
$ws.Range("C34").Value = "Spend"
$ws.Range("D34").Value = 127

$ws.Range("A35").Value = 46254
$ws.Range("A35").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B35").Value = "Radio"
$ws.Range("C35").Value = "Spend"
$ws.Range("D35").Value = 70

$ws.Range("A36").Value = 46240
$ws.Range("A36").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B36").Value = "Radio"
$ws.Range("C36").Value = "Spend"
$ws.Range("D36").Value = 54

$ws.Range("A37").Value = 46191
$ws.Range("A37").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B37").Value = "Radio"
$ws.Range("C37").Value = "GRPs"
$ws.Range("D37").Value = 10

$ws.Range("A38").Value = 46191
$ws.Range("A38").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B38").Value = "TV"
$ws.Range("C38").Value = "GRPs"
$ws.Range("D38").Value = 6

$ws.Range("A39").Value = 46212
$ws.Range("A39").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B39").Value = "TV"
$ws.Range("D39").Value = 1

$ws.Range("A40").Value = 46212
$ws.Range("A40").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B40").Value = "Radio"
$ws.Range("C40").Value = "GRPs"

$ws.Range("C41").Value = "Spend"
$ws.Range("D41").Value = 113

$ws.Range("C42").Value = "GRPs"

$ws.Range("A43").Value = 46184
$ws.Range("A43").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B43").Value = "Radio"
$ws.Range("C43").Value = "GRPs"
$ws.Range("D43").Value = 3

$ws.Range("A44").Value = 46205
$ws.Range("A44").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B44").Value = "TV"
$ws.Range("C44").Value = "GRPs"
$ws.Range("D44").Value = 1

$ws.Range("A45").Value = 46226
$ws.Range("A45").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B45").Value = "TV"
$ws.Range("C45").Value = "Spend"
$ws.Range("D45").Value = 200

$ws.Range("A46").Value = 46184
$ws.Range("A46").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B46").Value = "Radio"
$ws.Range("C46").Value = "GRPs"
$ws.Range("D46").Value = 3

$ws.Range("A47").Value = 46212
$ws.Range("A47").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B47").Value = "TV"
$ws.Range("C47").Value = "GRPs"
$ws.Range("D47").Value = 1
